# Auto-generated script to apply scheduled market-data refresh to Ultros_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1309.0571
$ws.Range("I28").Value = 1244.1034
$ws.Range("J28").Value = 1623
$ws.Range("K28").Value = 1244.1034
$ws.Range("L28").Value = 1623
$ws.Range("M28").Value = -759.1034
$ws.Range("N28").Value = -2593
$ws.Range("H86").Value = 2458
$ws.Range("I86").Value = 2249.6
$ws.Range("K86").Value = 2249.6
$ws.Range("M86").Value = -1126.6
$ws.Range("H89").Value = 2458
$ws.Range("I89").Value = 2249.6
$ws.Range("K89").Value = 11248
$ws.Range("M89").Value = -5632
$ws.Range("H96").Value = 277.5263
$ws.Range("I96").Value = 231.36363
$ws.Range("J96").Value = 341
$ws.Range("K96").Value = 694.09089
$ws.Range("L96").Value = 1023
$ws.Range("M96").Value = 678.90911
$ws.Range("N96").Value = -3769
$ws.Range("H100").Value = 5786.35
$ws.Range("I100").Value = 3045.5715
$ws.Range("J100").Value = 7262.154
$ws.Range("K100").Value = 3045.5715
$ws.Range("L100").Value = 7262.154
$ws.Range("M100").Value = -2504.5715
$ws.Range("N100").Value = -8344.154
$ws.Range("H101").Value = 2077.6155
$ws.Range("I101").Value = 1076.1111
$ws.Range("J101").Value = 4331
$ws.Range("K101").Value = 3228.3333
$ws.Range("L101").Value = 12993
$ws.Range("M101").Value = -1606.3333
$ws.Range("N101").Value = -16237
$ws.Range("H113").Value = 5532.486
$ws.Range("I113").Value = 4407.7646
$ws.Range("J113").Value = 6594.722
$ws.Range("K113").Value = 4407.7646
$ws.Range("L113").Value = 6594.722
$ws.Range("M113").Value = -1153.7646
$ws.Range("N113").Value = -13102.722
$ws.Range("H132").Value = 1357.9656
$ws.Range("I132").Value = 1266.0834
$ws.Range("K132").Value = 3798.2502
$ws.Range("M132").Value = -1268.2502
$ws.Range("H139").Value = 99999.91
$ws.Range("J139").Value = 99999.91
$ws.Range("L139").Value = 99999.91
$ws.Range("N139").Value = -110279.91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13528.882
$ws.Range("I32").Value = 8664.95
$ws.Range("K32").Value = 8664.95
$ws.Range("M32").Value = -8377.95
$ws.Range("H45").Value = 3233.6667
$ws.Range("I45").Value = 2300.5715
$ws.Range("K45").Value = 2300.5715
$ws.Range("M45").Value = -1923.5715
$ws.Range("H97").Value = 774.8182
$ws.Range("I97").Value = 847
$ws.Range("J97").Value = 608.8
$ws.Range("K97").Value = 847
$ws.Range("L97").Value = 608.8
$ws.Range("M97").Value = -351
$ws.Range("N97").Value = -1600.8
$ws.Range("H102").Value = 13334443
$ws.Range("I102").Value = 839.2381
$ws.Range("K102").Value = 839.2381
$ws.Range("M102").Value = 782.7619
$ws.Range("H110").Value = 1996.8064
$ws.Range("I110").Value = 2034.5518
$ws.Range("J110").Value = 1449.5
$ws.Range("K110").Value = 2034.5518
$ws.Range("L110").Value = 1449.5
$ws.Range("M110").Value = 10.44820000000004
$ws.Range("N110").Value = -5539.5
$ws.Range("H132").Value = 4806.467
$ws.Range("I132").Value = 4498.4194
$ws.Range("K132").Value = 13495.2582
$ws.Range("M132").Value = -10965.2582
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 6453469
$ws.Range("J94").Value = 25003022
$ws.Range("L94").Value = 25003022
$ws.Range("N94").Value = -25003924
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 896.9286
$ws.Range("I22").Value = 842.8571
$ws.Range("J22").Value = 951
$ws.Range("K22").Value = 842.8571
$ws.Range("L22").Value = 951
$ws.Range("M22").Value = -492.8570999999999
$ws.Range("N22").Value = -1651
$ws.Range("H58").Value = 2269.2856
$ws.Range("I58").Value = 977.5
$ws.Range("J58").Value = 5498.75
$ws.Range("K58").Value = 977.5
$ws.Range("L58").Value = 5498.75
$ws.Range("M58").Value = -774.5
$ws.Range("N58").Value = -5904.75
$ws.Range("H132").Value = 3874.6
$ws.Range("I132").Value = 3264.5881
$ws.Range("K132").Value = 9793.764299999999
$ws.Range("M132").Value = -7263.764299999999
$ws.Range("H134").Value = 2798.2917
$ws.Range("I134").Value = 1445.2368
$ws.Range("K134").Value = 4335.7104
$ws.Range("M134").Value = -1800.7104
$ws.Range("H136").Value = 2269.2856
$ws.Range("I136").Value = 977.5
$ws.Range("J136").Value = 5498.75
$ws.Range("K136").Value = 2932.5
$ws.Range("L136").Value = 16496.25
$ws.Range("M136").Value = -382.5
$ws.Range("N136").Value = -21596.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4820.2593
$ws.Range("I5").Value = 543.5625
$ws.Range("J5").Value = 11040.909
$ws.Range("K5").Value = 1630.6875
$ws.Range("L5").Value = 33122.727
$ws.Range("M5").Value = -1518.6875
$ws.Range("N5").Value = -33346.727
$ws.Range("H135").Value = 4820.2593
$ws.Range("I135").Value = 543.5625
$ws.Range("J135").Value = 11040.909
$ws.Range("K135").Value = 4892.0625
$ws.Range("L135").Value = 99368.181
$ws.Range("M135").Value = -2357.0625
$ws.Range("N135").Value = -104438.181
$ws.Range("H136").Value = 3275.5557
$ws.Range("J136").Value = 3996.8333
$ws.Range("L136").Value = 11990.4999
$ws.Range("N136").Value = -22190.4999
$ws.Range("H137").Value = 1500
$ws.Range("I137").Value = 1000
$ws.Range("K137").Value = 3000
$ws.Range("M137").Value = 2100
$ws.Range("H140").Value = 1814.1428
$ws.Range("J140").Value = 1566.3334
$ws.Range("L140").Value = 4699.0002
$ws.Range("N140").Value = -15059.0002
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 16725769
$ws.Range("I80").Value = 72689.375
$ws.Range("J80").Value = 83338090
$ws.Range("K80").Value = 72689.375
$ws.Range("L80").Value = 83338090
$ws.Range("M80").Value = -71691.375
$ws.Range("N80").Value = -83340086
$ws.Range("H83").Value = 16725769
$ws.Range("I83").Value = 72689.375
$ws.Range("J83").Value = 83338090
$ws.Range("K83").Value = 363446.875
$ws.Range("L83").Value = 416690450
$ws.Range("M83").Value = -358454.875
$ws.Range("N83").Value = -416700434
$ws.Range("H134").Value = 81721
$ws.Range("J134").Value = 81721
$ws.Range("L134").Value = 245163
$ws.Range("N134").Value = -250233
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1311.75
$ws.Range("I22").Value = 1327.6666
$ws.Range("J22").Value = 1264
$ws.Range("K22").Value = 1327.6666
$ws.Range("L22").Value = 1264
$ws.Range("M22").Value = -1032.6666
$ws.Range("N22").Value = -1854
$ws.Range("H27").Value = 1311.75
$ws.Range("I27").Value = 1327.6666
$ws.Range("J27").Value = 1264
$ws.Range("K27").Value = 1327.6666
$ws.Range("L27").Value = 1264
$ws.Range("M27").Value = -1220.6666
$ws.Range("N27").Value = -1478
$ws.Range("H40").Value = 6144.4053
$ws.Range("I40").Value = 6141.815
$ws.Range("J40").Value = 6151.4
$ws.Range("K40").Value = 6141.815
$ws.Range("L40").Value = 6151.4
$ws.Range("M40").Value = -6005.815
$ws.Range("N40").Value = -6423.4
$ws.Range("H55").Value = 2255.8
$ws.Range("I55").Value = 1987.3077
$ws.Range("J55").Value = 4001
$ws.Range("K55").Value = 1987.3077
$ws.Range("L55").Value = 4001
$ws.Range("M55").Value = -1814.3077
$ws.Range("N55").Value = -4347
$ws.Range("H58").Value = 10000
$ws.Range("I58").Value = 10000
$ws.Range("K58").Value = 10000
$ws.Range("M58").Value = -9740
$ws.Range("H61").Value = 2740.121
$ws.Range("I61").Value = 1777.4
$ws.Range("K61").Value = 1777.4
$ws.Range("M61").Value = -1575.4
$ws.Range("H68").Value = 5514.3447
$ws.Range("I68").Value = 5008.9473
$ws.Range("J68").Value = 6474.6
$ws.Range("K68").Value = 5008.9473
$ws.Range("L68").Value = 6474.6
$ws.Range("M68").Value = -4259.9473
$ws.Range("N68").Value = -7972.6
$ws.Range("H71").Value = 5514.3447
$ws.Range("I71").Value = 5008.9473
$ws.Range("J71").Value = 6474.6
$ws.Range("K71").Value = 25044.7365
$ws.Range("L71").Value = 32373
$ws.Range("M71").Value = -21300.7365
$ws.Range("N71").Value = -39861
$ws.Range("H74").Value = 82857.14
$ws.Range("J74").Value = 82857.14
$ws.Range("L74").Value = 82857.14
$ws.Range("N74").Value = -84853.14
$ws.Range("H77").Value = 82857.14
$ws.Range("J77").Value = 82857.14
$ws.Range("L77").Value = 248571.42
$ws.Range("N77").Value = -258555.42
$ws.Range("H113").Value = 2740.121
$ws.Range("I113").Value = 1777.4
$ws.Range("K113").Value = 1777.4
$ws.Range("M113").Value = 392.5999999999999
$ws.Range("H122").Value = 6034.05
$ws.Range("I122").Value = 4690.4614
$ws.Range("K122").Value = 14071.3842
$ws.Range("M122").Value = -11621.3842
$ws.Range("H132").Value = 4994.069
$ws.Range("I132").Value = 4122.9565
$ws.Range("J132").Value = 8333.333
$ws.Range("K132").Value = 12368.8695
$ws.Range("L132").Value = 24999.999
$ws.Range("M132").Value = -9838.8695
$ws.Range("N132").Value = -30059.999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 27781052
$ws.Range("J62").Value = 37040120
$ws.Range("L62").Value = 37040120
$ws.Range("N62").Value = -37041368
$ws.Range("H65").Value = 27781052
$ws.Range("J65").Value = 37040120
$ws.Range("L65").Value = 185200600
$ws.Range("N65").Value = -185206840
$ws.Range("H96").Value = 74135
$ws.Range("I96").Value = 102471.4
$ws.Range("J96").Value = 3294
$ws.Range("K96").Value = 102471.4
$ws.Range("L96").Value = 3294
$ws.Range("M96").Value = -101098.4
$ws.Range("N96").Value = -6040
$ws.Range("H100").Value = 710.8077
$ws.Range("I100").Value = 744.9545
$ws.Range("K100").Value = 1489.909
$ws.Range("M100").Value = -948.9090000000001
$ws.Range("H122").Value = 3234.2942
$ws.Range("I122").Value = 2679.6924
$ws.Range("J122").Value = 5036.75
$ws.Range("K122").Value = 8039.0772
$ws.Range("L122").Value = 15110.25
$ws.Range("M122").Value = -5589.0772
$ws.Range("N122").Value = -20010.25
$ws.Range("H126").Value = 3589.4
$ws.Range("I126").Value = 3432.6667
$ws.Range("K126").Value = 10298.0001
$ws.Range("M126").Value = -7828.000100000001
